# Regenerate save_data to use K (strikeouts) instead of Strike# column.
# This updates the "K" column (column G) values for each pitch-by-pitch
# game row (rows 2-29) on Sheet1 with the newly computed strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 1
    3  = 2
    4  = 3
    5  = 1
    6  = 2
    7  = 2
    8  = 5
    9  = 5
    10 = 2
    11 = 3
    12 = 5
    13 = 3
    14 = 5
    15 = 2
    16 = 3
    17 = 4
    18 = 2
    19 = 5
    20 = 4
    21 = 4
    22 = 3
    23 = 1
    24 = 4
    25 = 1
    26 = 0
    27 = 1
    28 = 4
    29 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
